$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 13407.186
$ws.Cells.Item(43, 10).Value = 15110.556
$ws.Cells.Item(43, 12).Value = 15110.556
$ws.Cells.Item(43, 14).Value = -15248.556

$ws.Cells.Item(69, 8).Value = 7484.933
$ws.Cells.Item(69, 9).Value = 2170.3333
$ws.Cells.Item(69, 10).Value = 8813.583000000001
$ws.Cells.Item(69, 11).Value = 6510.999899999999
$ws.Cells.Item(69, 12).Value = 26440.749
$ws.Cells.Item(69, 13).Value = -5636.999899999999
$ws.Cells.Item(69, 14).Value = -28188.749

$ws.Cells.Item(72, 8).Value = 7484.933
$ws.Cells.Item(72, 9).Value = 2170.3333
$ws.Cells.Item(72, 10).Value = 8813.583000000001
$ws.Cells.Item(72, 11).Value = 19532.9997
$ws.Cells.Item(72, 12).Value = 79322.247
$ws.Cells.Item(72, 13).Value = -15164.9997
$ws.Cells.Item(72, 14).Value = -88058.247

$ws.Cells.Item(88, 8).Value = 6725.6665
$ws.Cells.Item(88, 10).Value = 6870.8
$ws.Cells.Item(88, 12).Value = 6870.8
$ws.Cells.Item(88, 14).Value = -7682.8

$ws.Cells.Item(91, 8).Value = 6725.6665
$ws.Cells.Item(91, 10).Value = 6870.8
$ws.Cells.Item(91, 12).Value = 6870.8
$ws.Cells.Item(91, 14).Value = -9678.799999999999

$ws.Cells.Item(137, 8).Value = 30306184
$ws.Cells.Item(137, 9).Value = 58825732
$ws.Cells.Item(137, 11).Value = 176477196
$ws.Cells.Item(137, 13).Value = -176474646

$ws.Cells.Item(139, 8).Value = 69999
$ws.Cells.Item(139, 10).Value = 69999
$ws.Cells.Item(139, 12).Value = 69999
$ws.Cells.Item(139, 14).Value = -80279

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 23258894
$ws.Cells.Item(32, 9).Value = 23811986
$ws.Cells.Item(32, 11).Value = 23811986
$ws.Cells.Item(32, 13).Value = -23811699

$ws.Cells.Item(45, 8).Value = 3311.4707
$ws.Cells.Item(45, 9).Value = 2185.7144
$ws.Cells.Item(45, 11).Value = 2185.7144
$ws.Cells.Item(45, 13).Value = -1808.7144

$ws.Cells.Item(74, 8).Value = 11906242
$ws.Cells.Item(74, 9).Value = 11906242
$ws.Cells.Item(74, 11).Value = 11906242
$ws.Cells.Item(74, 13).Value = -11905368

$ws.Cells.Item(77, 8).Value = 11906242
$ws.Cells.Item(77, 9).Value = 11906242
$ws.Cells.Item(77, 11).Value = 59531210
$ws.Cells.Item(77, 13).Value = -59526842

$ws.Cells.Item(125, 8).Value = 65357.5
$ws.Cells.Item(125, 10).Value = 65357.5
$ws.Cells.Item(125, 12).Value = 65357.5
$ws.Cells.Item(125, 14).Value = -75197.5

$ws.Cells.Item(127, 8).Value = 57565
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 57565
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 57565
$ws.Cells.Item(127, 13).Value = $null
$ws.Cells.Item(127, 14).Value = -67485

$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1954
$ws.Cells.Item(107, 9).Value = 1954
$ws.Cells.Item(107, 11).Value = 1954
$ws.Cells.Item(107, 13).Value = -34

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 59689.316
$ws.Cells.Item(31, 10).Value = 110254
$ws.Cells.Item(31, 12).Value = 110254
$ws.Cells.Item(31, 14).Value = -110844

$ws.Cells.Item(34, 8).Value = 59689.316
$ws.Cells.Item(34, 10).Value = 110254
$ws.Cells.Item(34, 12).Value = 110254
$ws.Cells.Item(34, 14).Value = -110658

$ws.Cells.Item(58, 8).Value = 4143.32
$ws.Cells.Item(58, 9).Value = 1971.1765
$ws.Cells.Item(58, 10).Value = 8759.125
$ws.Cells.Item(58, 11).Value = 1971.1765
$ws.Cells.Item(58, 12).Value = 8759.125
$ws.Cells.Item(58, 13).Value = -1768.1765
$ws.Cells.Item(58, 14).Value = -9165.125

$ws.Cells.Item(86, 8).Value = 8889.727999999999
$ws.Cells.Item(86, 9).Value = 3898.5
$ws.Cells.Item(86, 11).Value = 3898.5
$ws.Cells.Item(86, 13).Value = -2775.5

$ws.Cells.Item(89, 8).Value = 8889.727999999999
$ws.Cells.Item(89, 9).Value = 3898.5
$ws.Cells.Item(89, 11).Value = 19492.5
$ws.Cells.Item(89, 13).Value = -13876.5

$ws.Cells.Item(105, 8).Value = 4039
$ws.Cells.Item(105, 9).Value = 4031.4285
$ws.Cells.Item(105, 10).Value = 4056.6667
$ws.Cells.Item(105, 11).Value = 4031.4285
$ws.Cells.Item(105, 12).Value = 4056.6667
$ws.Cells.Item(105, 13).Value = -2284.4285
$ws.Cells.Item(105, 14).Value = -7550.6667

$ws.Cells.Item(136, 8).Value = 4143.32
$ws.Cells.Item(136, 9).Value = 1971.1765
$ws.Cells.Item(136, 10).Value = 8759.125
$ws.Cells.Item(136, 11).Value = 5913.529500000001
$ws.Cells.Item(136, 12).Value = 26277.375
$ws.Cells.Item(136, 13).Value = -3363.529500000001
$ws.Cells.Item(136, 14).Value = -31377.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 33.125
$ws.Cells.Item(38, 9).Value = 49.57143
$ws.Cells.Item(38, 10).Value = 20.333334
$ws.Cells.Item(38, 11).Value = 148.71429
$ws.Cells.Item(38, 12).Value = 61.000002
$ws.Cells.Item(38, 13).Value = 198.28571
$ws.Cells.Item(38, 14).Value = -755.000002

$ws.Cells.Item(114, 8).Value = 1012.375
$ws.Cells.Item(114, 10).Value = 1084.8334
$ws.Cells.Item(114, 12).Value = 3254.5002
$ws.Cells.Item(114, 14).Value = -9762.5002

$ws.Cells.Item(137, 8).Value = 80241.46000000001
$ws.Cells.Item(137, 9).Value = 1775.75
$ws.Cells.Item(137, 11).Value = 5327.25
$ws.Cells.Item(137, 13).Value = -227.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(25, 8).Value = 6546
$ws.Cells.Item(25, 10).Value = 9769
$ws.Cells.Item(25, 12).Value = 9769
$ws.Cells.Item(25, 14).Value = -10827

$ws.Cells.Item(123, 8).Value = 38246
$ws.Cells.Item(123, 10).Value = 38246
$ws.Cells.Item(123, 12).Value = 38246
$ws.Cells.Item(123, 14).Value = -43146

$ws.Cells.Item(126, 8).Value = 4874.5
$ws.Cells.Item(126, 9).Value = 4333
$ws.Cells.Item(126, 10).Value = 6499
$ws.Cells.Item(126, 11).Value = 12999
$ws.Cells.Item(126, 12).Value = 19497
$ws.Cells.Item(126, 13).Value = -10529
$ws.Cells.Item(126, 14).Value = -24437

$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 7396.4
$ws.Cells.Item(32, 9).Value = 506.5
$ws.Cells.Item(32, 10).Value = 11989.667
$ws.Cells.Item(32, 11).Value = 506.5
$ws.Cells.Item(32, 12).Value = 11989.667
$ws.Cells.Item(32, 13).Value = -189.5
$ws.Cells.Item(32, 14).Value = -12623.667

$ws.Cells.Item(46, 8).Value = 1571.7222
$ws.Cells.Item(46, 10).Value = 1668.1875
$ws.Cells.Item(46, 12).Value = 1668.1875
$ws.Cells.Item(46, 14).Value = -2044.1875

$ws.Cells.Item(74, 8).Value = 52497.5
$ws.Cells.Item(74, 9).Value = 49995
$ws.Cells.Item(74, 11).Value = 49995
$ws.Cells.Item(74, 13).Value = -48997

$ws.Cells.Item(77, 8).Value = 52497.5
$ws.Cells.Item(77, 9).Value = 49995
$ws.Cells.Item(77, 11).Value = 149985
$ws.Cells.Item(77, 13).Value = -144993

$ws.Cells.Item(132, 8).Value = 10692.308
$ws.Cells.Item(132, 9).Value = 3165
$ws.Cells.Item(132, 10).Value = 12950.5
$ws.Cells.Item(132, 11).Value = 9495
$ws.Cells.Item(132, 12).Value = 38851.5
$ws.Cells.Item(132, 13).Value = -6965
$ws.Cells.Item(132, 14).Value = -43911.5

$ws.Cells.Item(136, 8).Value = 5267.0244
$ws.Cells.Item(136, 9).Value = 2082.5
$ws.Cells.Item(136, 10).Value = 10786.866
$ws.Cells.Item(136, 11).Value = 6247.5
$ws.Cells.Item(136, 12).Value = 32360.598
$ws.Cells.Item(136, 13).Value = -3697.5
$ws.Cells.Item(136, 14).Value = -37460.598

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).Value = $null

$ws.Cells.Item(70, 8).Value = 26249.5
$ws.Cells.Item(70, 10).Value = 19999.5
$ws.Cells.Item(70, 12).Value = 19999.5
$ws.Cells.Item(70, 14).Value = -20629.5

$ws.Cells.Item(73, 8).Value = 26249.5
$ws.Cells.Item(73, 10).Value = 19999.5
$ws.Cells.Item(73, 12).Value = 19999.5
$ws.Cells.Item(73, 14).Value = -22183.5

$ws.Cells.Item(82, 8).Value = 45000
$ws.Cells.Item(82, 9).Value = 45000
$ws.Cells.Item(82, 11).Value = 45000
$ws.Cells.Item(82, 13).Value = -44617

$ws.Cells.Item(85, 8).Value = 45000
$ws.Cells.Item(85, 9).Value = 45000
$ws.Cells.Item(85, 11).Value = 45000
$ws.Cells.Item(85, 13).Value = -43674

$ws.Cells.Item(132, 8).Value = 10166.111
$ws.Cells.Item(132, 9).Value = 2569.2856
$ws.Cells.Item(132, 11).Value = 7707.8568
$ws.Cells.Item(132, 13).Value = -5177.8568

$ws.Cells.Item(136, 8).Value = 1815.6578
$ws.Cells.Item(136, 9).Value = 1402.7576
$ws.Cells.Item(136, 11).Value = 4208.2728
$ws.Cells.Item(136, 13).Value = -1658.2728
